$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 7
    3  = -1
    6  = -5
    7  = 8
    8  = 1
    9  = -1
    10 = 0
    11 = 3
    12 = -3
    13 = 2
    15 = -2
    16 = 22
    17 = -1
    18 = 1
    19 = 7
    20 = -2
    21 = -5
    22 = -2
    23 = -4
    24 = -1
    25 = 7
    26 = -1
    28 = 7
    30 = 1
    31 = 3
    32 = -2
    33 = -1
    34 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
